# Auto-generated cell updates derived from the Adamantoise_Profits diff
# (market-data refresh from the scheduled Sheets runner)
$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3495.8
$ws.Range("I6").Value = 3495.8
$ws.Range("K6").Value = 10487.4
$ws.Range("M6").Value = -10375.4
$ws.Range("H17").Value = 8388.333000000001
$ws.Range("J17").Value = 8388.333000000001
$ws.Range("L17").Value = 25164.999
$ws.Range("N17").Value = -25500.999
$ws.Range("H21").Value = 21598.6
$ws.Range("I21").Value = 19999.5
$ws.Range("J21").Value = 22664.666
$ws.Range("K21").Value = 19999.5
$ws.Range("L21").Value = 22664.666
$ws.Range("M21").Value = -19531.5
$ws.Range("N21").Value = -23600.666
$ws.Range("H23").Value = 21598.6
$ws.Range("I23").Value = 19999.5
$ws.Range("J23").Value = 22664.666
$ws.Range("K23").Value = 19999.5
$ws.Range("L23").Value = 22664.666
$ws.Range("M23").Value = -19765.5
$ws.Range("N23").Value = -23132.666
$ws.Range("H28").Value = 1842.4286
$ws.Range("I28").Value = 1842.4286
$ws.Range("K28").Value = 1842.4286
$ws.Range("M28").Value = -1357.4286
$ws.Range("H31").Value = 996.6667
$ws.Range("J31").Value = 995
$ws.Range("L31").Value = 2985
$ws.Range("N31").Value = -3445
$ws.Range("H64").Value = 5142.923
$ws.Range("I64").Value = 4430.4287
$ws.Range("K64").Value = 4430.4287
$ws.Range("M64").Value = -4182.4287
$ws.Range("H67").Value = 5142.923
$ws.Range("I67").Value = 4430.4287
$ws.Range("K67").Value = 4430.4287
$ws.Range("M67").Value = -3572.4287
$ws.Range("H108").Value = 99594.75
$ws.Range("J108").Value = 99594.75
$ws.Range("L108").Value = 99594.75
$ws.Range("N108").Value = -107274.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2548.6428
$ws.Range("I32").Value = 2548.6428
$ws.Range("K32").Value = 2548.6428
$ws.Range("M32").Value = -2261.6428
$ws.Range("H45").Value = 3408.9583
$ws.Range("I45").Value = 3016.1765
$ws.Range("K45").Value = 3016.1765
$ws.Range("M45").Value = -2639.1765
$ws.Range("H61").Value = 3400.5356
$ws.Range("I61").Value = 3133.9583
$ws.Range("K61").Value = 3133.9583
$ws.Range("M61").Value = -2921.9583
$ws.Range("H63").Value = 159482.78
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 204335
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 204335
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -205707
$ws.Range("H66").Value = 159482.78
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 204335
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 1021675
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -1028539
$ws.Range("H74").Value = 1556.7
$ws.Range("I74").Value = 1481.6786
$ws.Range("J74").Value = 2607
$ws.Range("K74").Value = 1481.6786
$ws.Range("L74").Value = 2607
$ws.Range("M74").Value = -607.6786
$ws.Range("N74").Value = -4355
$ws.Range("H77").Value = 1556.7
$ws.Range("I77").Value = 1481.6786
$ws.Range("J77").Value = 2607
$ws.Range("K77").Value = 7408.393
$ws.Range("L77").Value = 13035
$ws.Range("M77").Value = -3040.393
$ws.Range("N77").Value = -21771
$ws.Range("H122").Value = 4402.375
$ws.Range("I122").Value = 3457.6775
$ws.Range("J122").Value = 7656.3335
$ws.Range("K122").Value = 10373.0325
$ws.Range("L122").Value = 22969.0005
$ws.Range("M122").Value = -7923.032499999999
$ws.Range("N122").Value = -27869.0005
$ws.Range("H131").Value = 64999
$ws.Range("J131").Value = 44999
$ws.Range("L131").Value = 44999
$ws.Range("N131").Value = -55079
$ws.Range("H136").Value = 3400.5356
$ws.Range("I136").Value = 3133.9583
$ws.Range("K136").Value = 9401.874899999999
$ws.Range("M136").Value = -6851.874899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 740.6923
$ws.Range("J80").Value = 752.875
$ws.Range("L80").Value = 752.875
$ws.Range("N80").Value = -2748.875
$ws.Range("H83").Value = 740.6923
$ws.Range("J83").Value = 752.875
$ws.Range("L83").Value = 3764.375
$ws.Range("N83").Value = -13748.375
$ws.Range("H86").Value = 3346.625
$ws.Range("I86").Value = 2419.9312
$ws.Range("K86").Value = 2419.9312
$ws.Range("M86").Value = -1296.9312
$ws.Range("H89").Value = 3346.625
$ws.Range("I89").Value = 2419.9312
$ws.Range("K89").Value = 12099.656
$ws.Range("M89").Value = -6483.655999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3567.4849
$ws.Range("I31").Value = 2135.5833
$ws.Range("K31").Value = 2135.5833
$ws.Range("M31").Value = -1840.5833
$ws.Range("H34").Value = 3567.4849
$ws.Range("I34").Value = 2135.5833
$ws.Range("K34").Value = 2135.5833
$ws.Range("M34").Value = -1933.5833
$ws.Range("H75").Value = 90419.28999999999
$ws.Range("J75").Value = 102787.4
$ws.Range("L75").Value = 102787.4
$ws.Range("N75").Value = -104783.4
$ws.Range("H78").Value = 90419.28999999999
$ws.Range("J78").Value = 102787.4
$ws.Range("L78").Value = 308362.2
$ws.Range("N78").Value = -318346.2
$ws.Range("H107").Value = 1548.7826
$ws.Range("I107").Value = 1284.7693
$ws.Range("K107").Value = 1284.7693
$ws.Range("M107").Value = 635.2307000000001
$ws.Range("H132").Value = 999.625
$ws.Range("I132").Value = 999.625
$ws.Range("K132").Value = 2998.875
$ws.Range("M132").Value = -468.875
$ws.Range("H134").Value = 1806.2046
$ws.Range("I134").Value = 1409.641
$ws.Range("K134").Value = 4228.923000000001
$ws.Range("M134").Value = -1693.923000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 79964140
$ws.Range("I4").Value = 71737070
$ws.Range("J4").Value = 93998570
$ws.Range("K4").Value = 215211210
$ws.Range("L4").Value = 281995710
$ws.Range("M4").Value = -215211098
$ws.Range("N4").Value = -281995934
$ws.Range("H12").Value = 22.916666
$ws.Range("I12").Value = 19.75
$ws.Range("J12").Value = 24.5
$ws.Range("K12").Value = 59.25
$ws.Range("L12").Value = 73.5
$ws.Range("M12").Value = 113.75
$ws.Range("N12").Value = -419.5
$ws.Range("H93").Value = 941.3333
$ws.Range("I93").Value = 824
$ws.Range("K93").Value = 2472
$ws.Range("M93").Value = -600
$ws.Range("H131").Value = 2335.375
$ws.Range("I131").Value = 1253.2
$ws.Range("K131").Value = 3759.6
$ws.Range("M131").Value = 1280.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3333
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 2999
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 2999
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -4995
$ws.Range("H83").Value = 3333
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 2999
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 14995
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -24979
$ws.Range("H95").Value = 39989
$ws.Range("J95").Value = 39989
$ws.Range("L95").Value = 39989
$ws.Range("N95").Value = -45481
$ws.Range("H102").Value = 2051.7646
$ws.Range("I102").Value = 1644.6154
$ws.Range("K102").Value = 1644.6154
$ws.Range("M102").Value = -22.61539999999991

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3428.875
$ws.Range("I46").Value = 461.375
$ws.Range("K46").Value = 461.375
$ws.Range("M46").Value = -273.375
$ws.Range("H101").Value = 50241
$ws.Range("J101").Value = 50241
$ws.Range("L101").Value = 50241
$ws.Range("N101").Value = -56731
$ws.Range("H122").Value = 2299.7
$ws.Range("I122").Value = 1966.1666
$ws.Range("K122").Value = 5898.4998
$ws.Range("M122").Value = -3448.4998
$ws.Range("H128").Value = 91750
$ws.Range("J128").Value = 91750
$ws.Range("L128").Value = 91750
$ws.Range("N128").Value = -101710
$ws.Range("H138").Value = 68999.5
$ws.Range("J138").Value = 68999.5
$ws.Range("L138").Value = 68999.5
$ws.Range("N138").Value = -79279.5
$ws.Range("H140").Value = 57404.168
$ws.Range("J140").Value = 41085
$ws.Range("L140").Value = 41085
$ws.Range("N140").Value = -51445
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 164452.75
$ws.Range("J128").Value = 164452.75
$ws.Range("L128").Value = 164452.75
$ws.Range("N128").Value = -174412.75
